# Realestate Update resale numbers 2023-06-24 09:58
# Appends the newest resale-number snapshot as row 72 of the
# "CityResaleNum" sheet (immediately after the existing last row, 71).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 72

# Date / Time / Weekday / Week are stored as plain text in this sheet
# (not real Excel dates/numbers), so force text entry with a leading
# apostrophe for the values that would otherwise be auto-detected as a
# date, time or number.
$ws.Range("A$row").Value = "'2023-06-24"
$ws.Range("B$row").Value = "'09:57:46"
$ws.Range("C$row").Value = "Saturday"
$ws.Range("D$row").Value = "'25"

# City resale counts - numeric.
$ws.Range("E$row").Value = 122564
$ws.Range("F$row").Value = 134108
$ws.Range("G$row").Value = 162588
$ws.Range("H$row").Value = 133495
$ws.Range("I$row").Value = 177463
$ws.Range("J$row").Value = 115576
$ws.Range("K$row").Value = 202709
$ws.Range("L$row").Value = 225828
$ws.Range("M$row").Value = 175769
$ws.Range("N$row").Value = 104145
$ws.Range("O$row").Value = 39455
$ws.Range("P$row").Value = 33818
$ws.Range("Q$row").Value = 51921
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 35743
$ws.Range("T$row").Value = -1
